$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1750
$ws.Range("F5").Value = 464
$ws.Range("F7").Value = 253
$ws.Range("F8").Value = 1214
$ws.Range("F9").Value = 349
$ws.Range("F12").Value = 696
$ws.Range("F13").Value = 190
$ws.Range("F14").Value = 518
$ws.Range("F15").Value = 143
$ws.Range("F18").Value = 2945
$ws.Range("F24").Value = 232
$ws.Range("F26").Value = 5309
$ws.Range("F29").Value = 24
$ws.Range("F31").Value = 319
$ws.Range("F32").Value = 1105
$ws.Range("F33").Value = 71
$ws.Range("F35").Value = 290
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1129
$ws.Range("F6").Value = 19
$ws.Range("F10").Value = 33
$ws.Range("F25").Value = 279
$ws.Range("F26").Value = 3950
$ws.Range("F33").Value = 167
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1783
$ws.Range("F5").Value = 2471
$ws.Range("F6").Value = 1051
$ws.Range("F9").Value = 1333
$ws.Range("F10").Value = 366
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1783
$ws.Range("F4").Value = 2471
$ws.Range("F5").Value = 1750
$ws.Range("F6").Value = 1051
$ws.Range("F7").Value = 1333
$ws.Range("F8").Value = 366
$ws.Range("F11").Value = 464
$ws.Range("F13").Value = 253
$ws.Range("F14").Value = 1214
$ws.Range("F15").Value = 349
$ws.Range("F17").Value = 696
$ws.Range("F18").Value = 1129
$ws.Range("F19").Value = 1129
$ws.Range("F20").Value = 190
$ws.Range("F21").Value = 518
$ws.Range("F22").Value = 19
$ws.Range("F23").Value = 2945
$ws.Range("F27").Value = 33
$ws.Range("F28").Value = 232
$ws.Range("F29").Value = 5309
$ws.Range("F34").Value = 24
$ws.Range("F37").Value = 319
$ws.Range("F45").Value = 279
$ws.Range("F46").Value = 1105
$ws.Range("F49").Value = 167
$ws.Range("F51").Value = 290
